# [1.10] ~ Generate single contact documents
# Add a new version-history row to the "Table2" ListObject on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the table by adding a new row at the bottom (this also grows the
# table/autoFilter range and the sheet dimension).
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Fill in the new row's values.
$ws.Range("A10").Value() = "[1.10]"
$ws.Range("B10").Value() = "~ Generate single contact documents"
$ws.Range("C10").Value() = 43258

# Match formatting of the existing rows: Version/Details columns use the
# wrapped, top-left aligned style (same as row 5), and the Date column uses
# the date-formatted style (same as row 8).
$ws.Range("A5:B5").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Update the selection/view to match the saved state of the workbook.
$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null
